# Update "想去人数" (want-to-go count) values in column F
# for worksheets "展览" and "全部类型".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 53
$ws1.Range("F6").Value = 423
$ws1.Range("F7").Value = 242
$ws1.Range("F8").Value = 13345
$ws1.Range("F10").Value = 50
$ws1.Range("F11").Value = 5403
$ws1.Range("F15").Value = 45
$ws1.Range("F16").Value = 1211
$ws1.Range("F17").Value = 57
$ws1.Range("F18").Value = 149
$ws1.Range("F19").Value = 716
$ws1.Range("F21").Value = 7390
$ws1.Range("F23").Value = 3667
$ws1.Range("F24").Value = 228

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 53
$ws4.Range("F7").Value = 423
$ws4.Range("F8").Value = 242
$ws4.Range("F9").Value = 13345
$ws4.Range("F11").Value = 50
$ws4.Range("F12").Value = 5403
$ws4.Range("F16").Value = 45
$ws4.Range("F17").Value = 1211
$ws4.Range("F18").Value = 57
$ws4.Range("F19").Value = 149
$ws4.Range("F20").Value = 716
$ws4.Range("F23").Value = 7390
$ws4.Range("F25").Value = 3667
$ws4.Range("F26").Value = 228

$wb.Save()
